$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 3175.199951171875
$ws.Range("E2").Value = 3187.8359375
$ws.Range("F2").Value = 3194.639892578125
$ws.Range("G2").Value = 3146.0400390625
$ws.Range("H2").Value = 5703702
$ws.Range("I2").Value = "KITT"
$ws.Range("D3").Value = 3175.199951171875
$ws.Range("E3").Value = 3187.8359375
$ws.Range("F3").Value = 3194.639892578125
$ws.Range("G3").Value = 3146.0400390625
$ws.Range("H3").Value = 5703702
$ws.Range("I3").Value = "KITT"
$ws.Range("D4").Value = 3175.199951171875
$ws.Range("E4").Value = 3187.8359375
$ws.Range("F4").Value = 3194.639892578125
$ws.Range("G4").Value = 3146.0400390625
$ws.Range("H4").Value = 5703702
$ws.Range("I4").Value = "KITT"
$ws.Range("D5").Value = 3175.199951171875
$ws.Range("E5").Value = 3187.8359375
$ws.Range("F5").Value = 3194.639892578125
$ws.Range("G5").Value = 3146.0400390625
$ws.Range("H5").Value = 5703702
$ws.Range("I5").Value = "KITT"
$ws.Range("D6").Value = 3217.320068359375
$ws.Range("E6").Value = 3223.800048828125
$ws.Range("F6").Value = 3541.320068359375
$ws.Range("G6").Value = 3210.840087890625
$ws.Range("H6").Value = 5703702
$ws.Range("I6").Value = "KITT"
$ws.Range("D7").Value = 3223.800048828125
$ws.Range("E7").Value = 3220.56005859375
$ws.Range("F7").Value = 3252.9599609375
$ws.Range("G7").Value = 3207.60009765625
$ws.Range("H7").Value = 5703702
$ws.Range("I7").Value = "KITT"
$ws.Range("D8").Value = 3256.199951171875
$ws.Range("E8").Value = 3256.199951171875
$ws.Range("F8").Value = 3262.679931640625
$ws.Range("G8").Value = 3249.719970703125
$ws.Range("H8").Value = 5703702
$ws.Range("I8").Value = "KITT"
$ws.Range("D9").Value = 3265.919921875
$ws.Range("E9").Value = 3298.320068359375
$ws.Range("F9").Value = 3434.39990234375
$ws.Range("G9").Value = 3262.679931640625
$ws.Range("H9").Value = 5703702
$ws.Range("I9").Value = "KITT"
$ws.Range("D10").Value = 1124.280029296875
$ws.Range("E10").Value = 1098.359985351562
$ws.Range("F10").Value = 2319.840087890625
$ws.Range("G10").Value = 1017.359985351562
$ws.Range("H10").Value = 5703702
$ws.Range("I10").Value = "KITT"
$ws.Range("D11").Value = 1215
$ws.Range("E11").Value = 1101.599975585938
$ws.Range("F11").Value = 1260.359985351562
$ws.Range("G11").Value = 1053
$ws.Range("H11").Value = 5703702
$ws.Range("I11").Value = "KITT"
$ws.Range("D12").Value = 972
$ws.Range("E12").Value = 761.4000244140625
$ws.Range("F12").Value = 1007.640014648438
$ws.Range("G12").Value = 725.760009765625
$ws.Range("H12").Value = 5703702
$ws.Range("I12").Value = "KITT"
$ws.Range("D13").Value = 657.719970703125
$ws.Range("E13").Value = 648
$ws.Range("F13").Value = 722.52001953125
$ws.Range("G13").Value = 622.0800170898438
$ws.Range("H13").Value = 5703702
$ws.Range("I13").Value = "KITT"
$ws.Range("D14").Value = 576.719970703125
$ws.Range("E14").Value = 589.6799926757812
$ws.Range("F14").Value = 638.280029296875
$ws.Range("G14").Value = 557.280029296875
$ws.Range("H14").Value = 5703702
$ws.Range("I14").Value = "KITT"
$ws.Range("D15").Value = 223.5599975585937
$ws.Range("E15").Value = 141.2640075683594
$ws.Range("F15").Value = 227.1239929199219
$ws.Range("G15").Value = 90.72000122070312
$ws.Range("H15").Value = 5703702
$ws.Range("I15").Value = "KITT"
$ws.Range("D16").Value = 105.3000030517578
$ws.Range("E16").Value = 67.39199829101562
$ws.Range("F16").Value = 105.3000030517578
$ws.Range("G16").Value = 59.29199981689453
$ws.Range("H16").Value = 5703702
$ws.Range("I16").Value = "KITT"
$ws.Range("D17").Value = 45.0359992980957
$ws.Range("E17").Value = 28.70999908447266
$ws.Range("F17").Value = 55.08000183105469
$ws.Range("G17").Value = 23.85000038146973
$ws.Range("H17").Value = 5703702
$ws.Range("I17").Value = "KITT"
$ws.Range("D18").Value = 12.86999988555908
$ws.Range("E18").Value = 12.32999992370606
$ws.Range("F18").Value = 13.5
$ws.Range("G18").Value = 10.89000034332275
$ws.Range("H18").Value = 5703702
$ws.Range("I18").Value = "KITT"
$ws.Range("D19").Value = 15.02999973297119
$ws.Range("E19").Value = 15.84000015258789
$ws.Range("F19").Value = 54.36000061035156
$ws.Range("G19").Value = 13.5
$ws.Range("H19").Value = 5703702
$ws.Range("I19").Value = "KITT"
$ws.Range("D20").Value = 8.1899995803833
$ws.Range("E20").Value = 9.09000015258789
$ws.Range("F20").Value = 12.51000022888184
$ws.Range("G20").Value = 7.380000114440918
$ws.Range("H20").Value = 5703702
$ws.Range("I20").Value = "KITT"
$ws.Range("D21").Value = 8.100000381469727
$ws.Range("E21").Value = 8.442000389099121
$ws.Range("F21").Value = 10.89000034332275
$ws.Range("G21").Value = 7.920000076293945
$ws.Range("H21").Value = 5703702
$ws.Range("I21").Value = "KITT"
